$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: tasks get reshuffled upward / new tasks appended at bottom ---
$ws.Range("B3").Value  = "توضیحات برند`nجدا از متن توضیحات محصول"
$ws.Range("B4").Value  = "لینک های نقد برررسی محصول"
$ws.Range("B5").Value  = "صفحه هر محصول انجمن"
$ws.Range("B6").Value  = "قیمت های اضافی ثابت برای محصولات مانند مالیات یا حمل و نقل"
$ws.Range("B7").Value  = "confirm email and phone number"
$ws.Range("B8").Value  = "ارسال اس ام اس"
$ws.Range("B9").Value  = "ارسال ایمیل"
$ws.Range("B10").Value = "زمانی که محصول به سبد خرید اضافه می شود. صفحه اتمام موجودی"
$ws.Range("B11").Value = "اطلاعات فاکتور ا و سبد خرید در ادمین"
$ws.Range("B12").Value = "کاربران و نقش ها در ادمین"
$ws.Range("B13").Value = "نهایی کردن سفارش"

# rows 14 and 15 no longer exist - clear their old contents
$ws.Range("B14").Clear()
$ws.Range("B15").Clear()

# --- Column C: "Doing" column populated with blog / product SEO items ---
$ws.Range("C2").Value = "لینک های پیشنهادی بلگ"
$ws.Range("C3").Value = "متا برای سئو بلگ"
$ws.Range("C4").Value = "مطالب مشابه بلاگ"
$ws.Range("C5").Value = "لینک های پیشنهادی محصول"
$ws.Range("C6").Value = "متا برای سئو محصول"

# --- Column D: "Check" column gets the old C2 task moved over ---
$ws.Range("D2").Value = "سفارشات"

# --- Column E: "Publish" column gets three more completed items ---
$ws.Range("E11").Value = "محصولات در دسته بندی محصول"
$ws.Range("E12").Value = "سبد خرید عملیاتی"
$ws.Range("E13").Value = "حذف از سبد خرید"

# new shared string order matches: "سئو در صفحه اصلی" is appended after "حذف از سبد خرید"
$ws.Range("C7").Value = "سئو در صفحه اصلی"

# --- Row heights (auto height from wrapped text, matches new content) ---
$ws.Rows.Item(2).RowHeight  = 60
$ws.Rows.Item(3).RowHeight  = 60
$ws.Rows.Item(4).RowHeight  = 60
$ws.Rows.Item(5).RowHeight  = 45
$ws.Rows.Item(6).RowHeight  = 90
$ws.Rows.Item(7).RowHeight  = 60
# Rows 8 and 9 now hold short single-line text, so they revert to the
# sheet's default (automatic) row height instead of their old custom one.
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(10).RowHeight = 90
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30

# rows 14/15 are gone now - clear any leftover row-specific formatting
$ws.Rows.Item(14).Clear()
$ws.Rows.Item(15).Clear()

# --- Selection moves to C8 ---
$ws.Range("C8").Select()
